$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A112").Value = "EXIT TO DESKTOP"
$ws.Range("B112").Value = "Exit to Desktop"

$ws.Range("B3").Value = "Danger Chasers"

$ws.Range("A113").Value = "KEYBINDINGS"
$ws.Range("B113").Value = "Keybindings"

$ws.Range("A114").Value = "PRESS ANY BUTTON"
$ws.Range("B114").Value = "Press Any Button"

$ws.Range("B114").Select() | Out-Null
